$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new "2021" column (L) by copying the existing "2020" column (K) so
# the new cells inherit the same styles/borders, then overwrite the values
# that differ for the new year.
$ws.Range("K3:K5").Copy($ws.Range("L3"))

$ws.Range("L4").Value = 2021
$ws.Range("L5").Value = 269

# Leave the selection where the user ended up after the edit.
$ws.Range("N3").Select()
